$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.117734
$ws.Range("H2").Value = 0.353202
$ws.Range("I2").Value = 0.6076081328197709
$ws.Range("J2").Value = 0.6076081328197709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.67700833333333
$ws.Range("N2").Value = 41.031025
$ws.Range("O2").Value = 0.124413831206147
$ws.Range("P2").Value = 0.124413831206147
$ws.Range("Q2").Value = 1.610248899116667
$ws.Range("R2").Value = 14.49224009205
$ws.Range("S2").Value = 0.0755948556761211
$ws.Range("T2").Value = 0.0755948556761211

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.117734
$ws.Range("H3").Value = 0.353202
$ws.Range("I3").Value = 0.6076081328197709
$ws.Range("J3").Value = 0.6076081328197709
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 74.64939600000001
$ws.Range("N3").Value = 223.948188
$ws.Range("O3").Value = 0.679053278848249
$ws.Range("P3").Value = 0.6790532788482488
$ws.Range("Q3").Value = 8.788771988664001
$ws.Range("R3").Value = 79.098947897976
$ws.Range("S3").Value = 0.4125982948461278
$ws.Range("T3").Value = 0.4125982948461277

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.117734
$ws.Range("H4").Value = 0.353202
$ws.Range("I4").Value = 0.6076081328197709
$ws.Range("J4").Value = 0.6076081328197709
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.629335666666667
$ws.Range("N4").Value = 4.888007
$ws.Range("O4").Value = 0.01482136207497777
$ws.Range("P4").Value = 0.01482136207497777
$ws.Range("Q4").Value = 0.1918282053793333
$ws.Range("R4").Value = 1.726453848414
$ws.Range("S4").Value = 0.00900558013622301
$ws.Range("T4").Value = 0.00900558013622301

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.117734
$ws.Range("H5").Value = 0.353202
$ws.Range("I5").Value = 0.6076081328197709
$ws.Range("J5").Value = 0.6076081328197709
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.17462033333333
$ws.Range("N5").Value = 57.523861
$ws.Range("O5").Value = 0.174423230537864
$ws.Range("P5").Value = 0.174423230537864
$ws.Range("Q5").Value = 2.257504750324667
$ws.Range("R5").Value = 20.317542752922
$ws.Range("S5").Value = 0.105980973427504
$ws.Range("T5").Value = 0.105980973427504

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.117734
$ws.Range("H6").Value = 0.353202
$ws.Range("I6").Value = 0.6076081328197709
$ws.Range("J6").Value = 0.6076081328197709
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.801214
$ws.Range("N6").Value = 2.403642
$ws.Range("O6").Value = 0.007288297332762355
$ws.Range("P6").Value = 0.007288297332762355
$ws.Range("Q6").Value = 0.094330129076
$ws.Range("R6").Value = 0.8489711616840001
$ws.Range("S6").Value = 0.004428428733795051
$ws.Range("T6").Value = 0.004428428733795051

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.029478
$ws.Range("H7").Value = 0.088434
$ws.Range("I7").Value = 0.1521316912638762
$ws.Range("J7").Value = 0.1521316912638762
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.67700833333333
$ws.Range("N7").Value = 41.031025
$ws.Range("O7").Value = 0.124413831206147
$ws.Range("P7").Value = 0.124413831206147
$ws.Range("Q7").Value = 0.40317085165
$ws.Range("R7").Value = 3.62853766485
$ws.Range("S7").Value = 0.01892728655800956
$ws.Range("T7").Value = 0.01892728655800956

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.029478
$ws.Range("H8").Value = 0.088434
$ws.Range("I8").Value = 0.1521316912638762
$ws.Range("J8").Value = 0.1521316912638762
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 74.64939600000001
$ws.Range("N8").Value = 223.948188
$ws.Range("O8").Value = 0.679053278848249
$ws.Range("P8").Value = 0.6790532788482488
$ws.Range("Q8").Value = 2.200514895288
$ws.Range("R8").Value = 19.804634057592
$ws.Range("S8").Value = 0.1033055237694647
$ws.Range("T8").Value = 0.1033055237694647

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.029478
$ws.Range("H9").Value = 0.088434
$ws.Range("I9").Value = 0.1521316912638762
$ws.Range("J9").Value = 0.1521316912638762
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.629335666666667
$ws.Range("N9").Value = 4.888007
$ws.Range("O9").Value = 0.01482136207497777
$ws.Range("P9").Value = 0.01482136207497777
$ws.Range("Q9").Value = 0.048029556782
$ws.Range("R9").Value = 0.432266011038
$ws.Range("S9").Value = 0.002254798879300643
$ws.Range("T9").Value = 0.002254798879300642

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Fgf15"
$ws.Range("C10").Value = "Fgfr1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.029478
$ws.Range("H10").Value = 0.088434
$ws.Range("I10").Value = 0.1521316912638762
$ws.Range("J10").Value = 0.1521316912638762
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.17462033333333
$ws.Range("N10").Value = 57.523861
$ws.Range("O10").Value = 0.174423230537864
$ws.Range("P10").Value = 0.174423230537864
$ws.Range("Q10").Value = 0.565229458186
$ws.Range("R10").Value = 5.087065123674
$ws.Range("S10").Value = 0.02653530105743424
$ws.Range("T10").Value = 0.02653530105743424

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Fgf15"
$ws.Range("C11").Value = "Fgfr1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.029478
$ws.Range("H11").Value = 0.088434
$ws.Range("I11").Value = 0.1521316912638762
$ws.Range("J11").Value = 0.1521316912638762
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.801214
$ws.Range("N11").Value = 2.403642
$ws.Range("O11").Value = 0.007288297332762355
$ws.Range("P11").Value = 0.007288297332762355
$ws.Range("Q11").Value = 0.023618186292
$ws.Range("R11").Value = 0.212563676628
$ws.Range("S11").Value = 0.001108780999667135
$ws.Range("T11").Value = 0.001108780999667135

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgf15"
$ws.Range("C12").Value = "Fgfr1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.04655433333333334
$ws.Range("H12").Value = 0.139663
$ws.Range("I12").Value = 0.2402601759163528
$ws.Range("J12").Value = 0.2402601759163529
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 13.67700833333333
$ws.Range("N12").Value = 41.031025
$ws.Range("O12").Value = 0.124413831206147
$ws.Range("P12").Value = 0.124413831206147
$ws.Range("Q12").Value = 0.6367240049527778
$ws.Range("R12").Value = 5.730516044575
$ws.Range("S12").Value = 0.0298916889720163
$ws.Range("T12").Value = 0.0298916889720163

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgf15"
$ws.Range("C13").Value = "Fgfr1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.04655433333333334
$ws.Range("H13").Value = 0.139663
$ws.Range("I13").Value = 0.2402601759163528
$ws.Range("J13").Value = 0.2402601759163529
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 74.64939600000001
$ws.Range("N13").Value = 223.948188
$ws.Range("O13").Value = 0.679053278848249
$ws.Range("P13").Value = 0.6790532788482488
$ws.Range("Q13").Value = 3.475252864516001
$ws.Range("R13").Value = 31.277275780644
$ws.Range("S13").Value = 0.1631494602326565
$ws.Range("T13").Value = 0.1631494602326565

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Fgf15"
$ws.Range("C14").Value = "Fgfr1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04655433333333334
$ws.Range("H14").Value = 0.139663
$ws.Range("I14").Value = 0.2402601759163528
$ws.Range("J14").Value = 0.2402601759163529
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.629335666666667
$ws.Range("N14").Value = 4.888007
$ws.Range("O14").Value = 0.01482136207497777
$ws.Range("P14").Value = 0.01482136207497777
$ws.Range("Q14").Value = 0.07585263573788889
$ws.Range("R14").Value = 0.682673721641
$ws.Range("S14").Value = 0.00356098305945412
$ws.Range("T14").Value = 0.00356098305945412

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Fgf15"
$ws.Range("C15").Value = "Fgfr1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04655433333333334
$ws.Range("H15").Value = 0.139663
$ws.Range("I15").Value = 0.2402601759163528
$ws.Range("J15").Value = 0.2402601759163529
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 19.17462033333333
$ws.Range("N15").Value = 57.523861
$ws.Range("O15").Value = 0.174423230537864
$ws.Range("P15").Value = 0.174423230537864
$ws.Range("Q15").Value = 0.8926616665381112
$ws.Range("R15").Value = 8.033954998843
$ws.Range("S15").Value = 0.04190695605292577
$ws.Range("T15").Value = 0.04190695605292578

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Fgf15"
$ws.Range("C16").Value = "Fgfr1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04655433333333334
$ws.Range("H16").Value = 0.139663
$ws.Range("I16").Value = 0.2402601759163528
$ws.Range("J16").Value = 0.2402601759163529
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.801214
$ws.Range("N16").Value = 2.403642
$ws.Range("O16").Value = 0.007288297332762355
$ws.Range("P16").Value = 0.007288297332762355
$ws.Range("Q16").Value = 0.03729998362733333
$ws.Range("R16").Value = 0.3356998526460001
$ws.Range("S16").Value = 0.001751087599300169
$ws.Range("T16").Value = 0.001751087599300169

